# Simulated Wild Card round and logged it
# Appends the new game's per-play data to the season logs (YDS, ST sheets)
# and updates the aggregated season totals (OFF, DEF, ST, TURNS, PEN sheets).

$wb = $excel.ActiveWorkbook

# ---- YDS sheet: per-play yardage logs (Home/Road x OFF/DEF) ----
$ydsSheet = $wb.Worksheets.Item("YDS")
$ydsSheet.Range("B2").Value = "-1 2 5 2 9 2 2 -2 8 7 3 0 5 4 6 6 2 2 12 3 8 0 2 4 70 3 -1 2 0 5 9 3 2 7 0 1 -6 0 3 3 16 2 2 5 0 5 8 -1 3 14 2 0 3 6 5 7 3 8 4 11 0 1 14 2 2 0 13 4 13 1 10 0 -1 3 2 -3 10 4 7 2 -1 1 -2 1 0 -2 9 2 2 6 -1 7 8 1 4 14 3 2 31 3 3 3 3 -3 6 10 -1 12 1 5 2 -4 8 0 49 2 7 1 1 2 3 10 -1 11 4 5 10 1 9 3 11 3 4 0 3 0 5 8 4 0 30 4 2 4 4 1 1 0 6 2 1 6 0 2 7 3 3 5 2 0 7 2 8 4 2 6 2 2 7 1 1 7 2 -1 5 3 0 11 14 1 7 1 -2 -2 11 20 -3 7 13 1 -2 3 17 0 12 2 2 11 2 2 2 30 2 1 1 7 14 6 7 6 4 -4 2 0 1 10 3 2 0 2 1 -1 1 3 9 11 2 15 7 1 1 4 20 -4 8 9 1 2 12 2 10 2 4 2 4 2 -2 11 12 8 16 -1 9 1 2 2 6 6 2 2 0 3 0 7 4 3 8 4 4 3 9 -1 1 14 2 9 6 7 2 14 3 0 1 11 3 5 5 -1 5 4 8 9 7 2 6 2 16 1 11 5 3 1 13 9 1 5 6 3 3 9 4 1 3 2 2 1 3 3 -1 6 16 5 5 5 9 14 -1 1 1 6 8 3 4 -1 10 5 1 4 2 4 4 3 5 -1 5 2 2 1 3 2 6 7 2 2 5 6 14 4 5 7 1 0 -1 7 7 13 1 -2 10 1 2 4 5 1 1 1 3 1 3 -4 1 -3 1 3 0 -4 5 3 6 7 -3 5 5 4 5 2 0 1 0 6 12 4 7 4 -4 -4 2 2 0 4 5 2 5 6 2 6 1 9 6 1 5 29 6 5 3 0 3 23 -1 4 3 -3 47 7 0 5 15 8 -4 4 1"
$ydsSheet.Range("C2").Value = "5 -2 4 2 -1 1 5 2 3 7 3 5 2 4 2 6 1 4 4 4 2 0 6 2 2 1 6 4 -1 12 7 11 5 4 1 4 2 2 2 13 8 5 4 2 -7 3 -1 7 5 2 1 1 3 4 5 7 11 2 1 -1 6 -2 3 4 3 5 4 0 2 3 5 3 5 5 4 5 6 3 3 3 5 7 6 5 2 4 5 20 -3 3 5 3 0 0 1 3 13 -2 5 18 3 8 -1 -1 4 9 2 1 4 4 1 4 6 2 6 4 3 9 3 6 18 6 2 0 1 3 1 5 3 1 3 0 20 10 4 12 1 9 2 1 2 5 10 11 2 3 3 2 0 9 1 6 1 9 5 -2 20 6 3 7 3 5 8 0 2 10 4 3 2 2 5 3 6 5 0 4 4 7 0 4 5 13 12 12 4 4 4 5 8 1 7 11 7 2 3 -2 7 2 5 10 9 10 17 -2 0 1 4 3 -8 7 4 2 5 11 31 0 3 4 3 0 5 -1 5 21 -2 15 5 25 1 7 5 3 1 2 6 1 23 1 1 7 5 6 3 3 5 1 5 -1 2 13 5 2 1 13 10 2 3 5 10 4 2 12 4 2 0 5 4 2 4 4 4 4 6 7 -4 1 5 1 13 35 5 9 3 4 4 2 3 5 4 -3 7 3 3 6 -3 8 6 0 2 4 5 5 6 5 1 -1 5 8 4 3 5 22 2 12 1 8 2 4 1 1 3 3 3 7 4 6 5 3 3 -2 0 -4 8 7 3 8 19 4 8 5 5 4 8 -4 3 0 4 4 4 10 6 0 4 13 3 5 4 9 -1 2 14 3 2 6 4 3 1 6 9 6 0 2 0 -1 1 2 -4 2 7 4 5 -4 2 5 7 2 4 1 3 0 11 3 3 3 11 3 0 1 6 10 5 6 5 25 5 3 6 7 4 6 4 4 2 5"
$ydsSheet.Range("B3").Value = "15 13 7 6 1 9 5 25 -4 5 11 6 2 6 17 16 2 10 6 8 20 7 22 13 3 6 13 14 8 6 6 7 4 10 16 33 15 12 9 55 14 5 10 11 9 0 4 11 14 6 6 10 8 5 28 6 17 6 9 27 6 21 -3 9 7 24 14 4 31 11 19 3 3 12 14 3 11 11 8 3 5 6 13 -3 9 9 32 9 5 10 9 6 14 12 12 3 12 19 10 9 15 9 39 23 6 9 5 23 6 7 14 23 5 1 3 6 11 7 5 7 10 14 3 14 11 6 10 5 7 26 13 12 12 12 26 6 18 16 6 4 4 1 6 6 8 31 1 0 6 6 9 12 13 4 11 9 8 8 5 2 9 10 3 7 9 10 1 31 5 31 12 19 19 15 7 3 4 -4 9 2 11 14 9 18 7 -3 25 5 8 11 44 19 7 6 9 17 6 3 11 20 8 -3 16 40 -2 7 5 64 4 13 15 26 10 7 4 10 6 15 7 5 -9 0 13 9 14 3 3 13 19 5 7 12 7 5 2 12 1 12 14 42 1 36 10 10 6 3 2 34 5 9 5 8 11 9 13 19 12 6 19 10 7 13 32 8 9 8 6 7 7 4 8 7 9 10 16 5 5 4 5 8 7 28 -1 10 3 11 -1 8 9 11 8 25 2 24 24 5 5 2 7 1 40 4 14 6 8 20 5 13 4 8 15 0 44 18 -1 24 8 35 7 9 5 11 6 2 16 7 34 5 6 16 11 8 2 31 11 6 13 7 10 29 28"
$ydsSheet.Range("C3").Value = "42 -6 4 17 1 6 8 37 2 15 20 8 17 14 6 16 16 6 9 15 8 6 4 24 7 13 25 -2 -8 6 5 19 0 6 6 13 11 4 8 7 6 5 6 -5 14 5 3 22 7 4 10 23 4 13 13 20 6 12 24 12 9 49 32 7 13 14 4 8 13 7 8 17 4 7 14 32 4 23 50 4 5 13 6 59 5 14 14 3 18 1 22 16 9 48 2 3 25 11 16 6 8 8 10 33 31 29 3 40 8 51 7 34 4 5 14 8 4 6 9 15 6 12 -2 6 17 16 18 1 4 8 7 8 27 17 2 13 -2 20 4 12 1 17 14 7 5 14 4 23 20 19 5 22 14 3 1 5 0 32 17 21 4 18 5 4 24 9 32 5 10 7 4 20 11 4 16 18 3 4 6 24 36 8 6 4 5 9 11 12 14 2 7 24 -2 10 5 9 13 6 18 14 12 9 11 13 17 8 8 4 2 7 13 18 34 9 5 -2 16 9 9 12 13 9 29 14 -7 12 6 14 14 38 8 4 1 13 7 7 2 19 7 5 20 11 -2 9 9 11 13 4 6 14 8 14 7 2 6 2 6 17 2 2 9 8 4 19 3 56 4 4 15 8 5 16 18 6 4 0 15 10 20 -5 5 7 13 12 12 18 6 3 28 10 8 10 11 8 3 6 8 3 5 9 3 8 7 40 16 5 5 45 8 8 3 21 11 4 9 4 3 3 11 11 10 11 2 11 17 10 7 14 5 10 4 1 19 4 44 8 12"

# ---- OFF sheet: offensive season totals ----
$offSheet = $wb.Worksheets.Item("OFF")
$offSheet.Range("C2").Value = 209
$offSheet.Range("F2").Value = 69
$offSheet.Range("G2").Value = 64
$offSheet.Range("I2").Value = 10
$offSheet.Range("J2").Value = 32
$offSheet.Range("N2").Value = 18
$offSheet.Range("B3").Value = 11
$offSheet.Range("C3").Value = 171
$offSheet.Range("E3").Value = 37
$offSheet.Range("F3").Value = 96
$offSheet.Range("G3").Value = 30
$offSheet.Range("H3").Value = 37
$offSheet.Range("I3").Value = 73
$offSheet.Range("J3").Value = 44
$offSheet.Range("L3").Value = 287
$offSheet.Range("M3").Value = 189
$offSheet.Range("Q3").Value = 544

# ---- DEF sheet: defensive season totals ----
$defSheet = $wb.Worksheets.Item("DEF")
$defSheet.Range("C2").Value = 209
$defSheet.Range("E2").Value = 10
$defSheet.Range("F2").Value = 70
$defSheet.Range("G2").Value = 53
$defSheet.Range("I2").Value = 8
$defSheet.Range("J2").Value = 40
$defSheet.Range("N2").Value = 11
$defSheet.Range("O2").Value = 26
$defSheet.Range("P2").Value = 11
$defSheet.Range("C3").Value = 160
$defSheet.Range("E3").Value = 32
$defSheet.Range("F3").Value = 99
$defSheet.Range("G3").Value = 51
$defSheet.Range("H3").Value = 29
$defSheet.Range("I3").Value = 56
$defSheet.Range("J3").Value = 67
$defSheet.Range("L3").Value = 296
$defSheet.Range("M3").Value = 176
$defSheet.Range("Q3").Value = 552

# ---- ST sheet: special-teams season totals + per-kick logs ----
$stSheet = $wb.Worksheets.Item("ST")
$stSheet.Range("B2").Value = 77
$stSheet.Range("D2").Value = 65
$stSheet.Range("F2").Value = 259
$stSheet.Range("G2").Value = 252
$stSheet.Range("J2").Value = 118
$stSheet.Range("K2").Value = 112
$stSheet.Range("B3").Value = 63
$stSheet.Range("B6").Value = "17 19 22 18 13 18 12 23 17 18 21 18 7 26 18 18 0"
$stSheet.Range("D3").Value = "37 41 41 45 48 46 43 53 26 50 47 57 49 39 56 44 28 68 53 52 52 39 33 47 33 46 38 46 59 47 50 44 48 43 41 56 60 44 52 42 43 40 56 56 48 47 41 42 51 57 44 48 63 46 54 28 56 40 45 41 56 40 65 44 36"
$stSheet.Range("D4").Value = "0 2 0 2 0 0 0 0 0 12 0 0 10 0 42 0 0 0 0 0 0 0 0 0 5 0 0 3 10 0 0 12 0 0 0 0 10 0 9 0 0 0 10 3 0 0 0 -4 12 0 0 5 0 5 0 0 6 0 0 0 1 4 5 0 0"
$stSheet.Range("D5").Value = "0 0 2 11 0 0 0 0 0 12 18 0 0 0 4 0 0"

# ---- TURNS sheet: turnover season totals ----
$turnsSheet = $wb.Worksheets.Item("TURNS")
$turnsSheet.Range("D3").Value = 8

# ---- PEN sheet: penalty season totals ----
$penSheet = $wb.Worksheets.Item("PEN")
$penSheet.Range("B2").Value = 18
$penSheet.Range("D2").Value = 6
$penSheet.Range("D4").Value = 6

Write-Output "Wild Card round simulated and logged."
